$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 17.11000061035156
$ws.Range("E2").Value = 14.59000015258789
$ws.Range("F2").Value = 17.97999954223633
$ws.Range("G2").Value = 14.38000011444092
$ws.Range("H2").Value = 112123800
$ws.Range("I2").Value = "TSEM"

$ws.Range("D3").Value = 15.72999954223633
$ws.Range("E3").Value = 13.26000022888184
$ws.Range("F3").Value = 15.77000045776367
$ws.Range("G3").Value = 12.84000015258789
$ws.Range("H3").Value = 112123800
$ws.Range("I3").Value = "TSEM"

$ws.Range("D4").Value = 12.92000007629394
$ws.Range("E4").Value = 13.35999965667725
$ws.Range("F4").Value = 14.19999980926514
$ws.Range("G4").Value = 12.03999996185303
$ws.Range("H4").Value = 112123800
$ws.Range("I4").Value = "TSEM"

$ws.Range("D5").Value = 13.80000019073486
$ws.Range("E5").Value = 12.56999969482422
$ws.Range("F5").Value = 14.5
$ws.Range("G5").Value = 10.35999965667725
$ws.Range("H5").Value = 112123800
$ws.Range("I5").Value = "TSEM"

$ws.Range("D6").Value = 12.01000022888184
$ws.Range("E6").Value = 11.6899995803833
$ws.Range("F6").Value = 12.52000045776367
$ws.Range("G6").Value = 11.0600004196167
$ws.Range("H6").Value = 112123800
$ws.Range("I6").Value = "TSEM"

$ws.Range("D7").Value = 12.30000019073486
$ws.Range("E7").Value = 13.47999954223633
$ws.Range("F7").Value = 13.71000003814697
$ws.Range("G7").Value = 11.73999977111816
$ws.Range("H7").Value = 112123800
$ws.Range("I7").Value = "TSEM"

$ws.Range("D8").Value = 15.19999980926514
$ws.Range("E8").Value = 15.55000019073486
$ws.Range("F8").Value = 16.47999954223633
$ws.Range("G8").Value = 15.03999996185303
$ws.Range("H8").Value = 112123800
$ws.Range("I8").Value = "TSEM"

$ws.Range("D9").Value = 19.25
$ws.Range("E9").Value = 21.20999908447266
$ws.Range("F9").Value = 21.5
$ws.Range("G9").Value = 19.02000045776367
$ws.Range("H9").Value = 112123800
$ws.Range("I9").Value = "TSEM"

$ws.Range("D10").Value = 23.15999984741211
$ws.Range("E10").Value = 21.52000045776367
$ws.Range("F10").Value = 23.20000076293945
$ws.Range("G10").Value = 20.60000038146973
$ws.Range("H10").Value = 112123800
$ws.Range("I10").Value = "TSEM"

$ws.Range("D11").Value = 24.28000068664551
$ws.Range("E11").Value = 26.14999961853028
$ws.Range("F11").Value = 27.60000038146973
$ws.Range("G11").Value = 23.3799991607666
$ws.Range("H11").Value = 112123800
$ws.Range("I11").Value = "TSEM"

$ws.Range("D12").Value = 31.29000091552734
$ws.Range("E12").Value = 33.06000137329102
$ws.Range("F12").Value = 33.29999923706055
$ws.Range("G12").Value = 30.39999961853028
$ws.Range("H12").Value = 112123800
$ws.Range("I12").Value = "TSEM"

$ws.Range("D13").Value = 34.54999923706055
$ws.Range("E13").Value = 34.7400016784668
$ws.Range("F13").Value = 36.08000183105469
$ws.Range("G13").Value = 33.0099983215332
$ws.Range("H13").Value = 112123800
$ws.Range("I13").Value = "TSEM"

$ws.Range("D14").Value = 26.64999961853028
$ws.Range("E14").Value = 25.8700008392334
$ws.Range("F14").Value = 29.36000061035156
$ws.Range("G14").Value = 25.03000068664551
$ws.Range("H14").Value = 112123800
$ws.Range("I14").Value = "TSEM"

$ws.Range("D15").Value = 21.55999946594238
$ws.Range("E15").Value = 20.64999961853028
$ws.Range("F15").Value = 23.54999923706055
$ws.Range("G15").Value = 19.61000061035156
$ws.Range("H15").Value = 112123800
$ws.Range("I15").Value = "TSEM"

$ws.Range("D16").Value = 21.79000091552734
$ws.Range("E16").Value = 15.6899995803833
$ws.Range("F16").Value = 22.40999984741211
$ws.Range("G16").Value = 13.86999988555908
$ws.Range("H16").Value = 112123800
$ws.Range("I16").Value = "TSEM"

$ws.Range("D17").Value = 14.36999988555908
$ws.Range("E17").Value = 14.89000034332275
$ws.Range("F17").Value = 15.1899995803833
$ws.Range("G17").Value = 13.56999969482422
$ws.Range("H17").Value = 112123800
$ws.Range("I17").Value = "TSEM"

$ws.Range("D18").Value = 16.76000022888184
$ws.Range("E18").Value = 17.96999931335449
$ws.Range("F18").Value = 18.28000068664551
$ws.Range("G18").Value = 16.76000022888184
$ws.Range("H18").Value = 112123800
$ws.Range("I18").Value = "TSEM"

$ws.Range("D19").Value = 16.64999961853027
$ws.Range("E19").Value = 19.67000007629395
$ws.Range("F19").Value = 20.45000076293945
$ws.Range("G19").Value = 15.88000011444092
$ws.Range("H19").Value = 112123800
$ws.Range("I19").Value = "TSEM"

$ws.Range("D20").Value = 19.40999984741211
$ws.Range("E20").Value = 21.93000030517578
$ws.Range("F20").Value = 22.84000015258789
$ws.Range("G20").Value = 19.34000015258789
$ws.Range("H20").Value = 112123800
$ws.Range("I20").Value = "TSEM"

$ws.Range("D21").Value = 24.31999969482422
$ws.Range("E21").Value = 22.18000030517578
$ws.Range("F21").Value = 25.79999923706055
$ws.Range("G21").Value = 22.09000015258789
$ws.Range("H21").Value = 112123800
$ws.Range("I21").Value = "TSEM"

$ws.Range("D22").Value = 15.48999977111816
$ws.Range("E22").Value = 19.23999977111816
$ws.Range("F22").Value = 19.86000061035156
$ws.Range("G22").Value = 15.07999992370606
$ws.Range("H22").Value = 112123800
$ws.Range("I22").Value = "TSEM"

$ws.Range("D23").Value = 19.03000068664551
$ws.Range("E23").Value = 21.5
$ws.Range("F23").Value = 23.84000015258789
$ws.Range("G23").Value = 18.88999938964844
$ws.Range("H23").Value = 112123800
$ws.Range("I23").Value = "TSEM"

$ws.Range("D24").Value = 18.38999938964844
$ws.Range("E24").Value = 21.10000038146973
$ws.Range("F24").Value = 21.38999938964844
$ws.Range("G24").Value = 18.20999908447266
$ws.Range("H24").Value = 112123800
$ws.Range("I24").Value = "TSEM"

$ws.Range("D25").Value = 25.69000053405762
$ws.Range("E25").Value = 27.96999931335449
$ws.Range("F25").Value = 30.8700008392334
$ws.Range("G25").Value = 25.69000053405762
$ws.Range("H25").Value = 112123800
$ws.Range("I25").Value = "TSEM"

$ws.Range("D26").Value = 28.5
$ws.Range("E26").Value = 28.29999923706055
$ws.Range("F26").Value = 30.60000038146973
$ws.Range("G26").Value = 27.89999961853028
$ws.Range("H26").Value = 112123800
$ws.Range("I26").Value = "TSEM"

$ws.Range("D27").Value = 29.39999961853028
$ws.Range("E27").Value = 27.85000038146973
$ws.Range("F27").Value = 29.76000022888184
$ws.Range("G27").Value = 25.54999923706055
$ws.Range("H27").Value = 112123800
$ws.Range("I27").Value = "TSEM"

$ws.Range("D28").Value = 30
$ws.Range("E28").Value = 31.8700008392334
$ws.Range("F28").Value = 32.65999984741211
$ws.Range("G28").Value = 28.42000007629395
$ws.Range("H28").Value = 112123800
$ws.Range("I28").Value = "TSEM"

$ws.Range("D29").Value = 40.20000076293945
$ws.Range("E29").Value = 34.34999847412109
$ws.Range("F29").Value = 41.31000137329102
$ws.Range("G29").Value = 31.04000091552734
$ws.Range("H29").Value = 112123800
$ws.Range("I29").Value = "TSEM"

$ws.Range("D30").Value = 48.5099983215332
$ws.Range("E30").Value = 48.29999923706055
$ws.Range("F30").Value = 49.0099983215332
$ws.Range("G30").Value = 48.09999847412109
$ws.Range("H30").Value = 112123800
$ws.Range("I30").Value = "TSEM"

$ws.Range("D31").Value = 46.06000137329102
$ws.Range("E31").Value = 47.86000061035156
$ws.Range("F31").Value = 47.97999954223633
$ws.Range("G31").Value = 45.34000015258789
$ws.Range("H31").Value = 112123800
$ws.Range("I31").Value = "TSEM"

$ws.Range("D32").Value = 43.79999923706055
$ws.Range("E32").Value = 42.7599983215332
$ws.Range("F32").Value = 45
$ws.Range("G32").Value = 42.04000091552734
$ws.Range("H32").Value = 112123800
$ws.Range("I32").Value = "TSEM"

$ws.Range("D33").Value = 43.34999847412109
$ws.Range("E33").Value = 41.79000091552734
$ws.Range("F33").Value = 44.40999984741211
$ws.Range("G33").Value = 41.72999954223633
$ws.Range("H33").Value = 112123800
$ws.Range("I33").Value = "TSEM"

$ws.Range("D34").Value = 42.4900016784668
$ws.Range("E34").Value = 45.0099983215332
$ws.Range("F34").Value = 45.86999893188477
$ws.Range("G34").Value = 41.79000091552734
$ws.Range("H34").Value = 112123800
$ws.Range("I34").Value = "TSEM"

$ws.Range("D35").Value = 37.25
$ws.Range("E35").Value = 37.75
$ws.Range("F35").Value = 39
$ws.Range("G35").Value = 34.7400016784668
$ws.Range("H35").Value = 112123800
$ws.Range("I35").Value = "TSEM"

$ws.Range("D36").Value = 24.64999961853028
$ws.Range("E36").Value = 23.02000045776367
$ws.Range("F36").Value = 24.84000015258789
$ws.Range("G36").Value = 21.43000030517578
$ws.Range("H36").Value = 112123800
$ws.Range("I36").Value = "TSEM"

$ws.Range("D37").Value = 30.53000068664551
$ws.Range("E37").Value = 28.85000038146973
$ws.Range("F37").Value = 30.82999992370605
$ws.Range("G37").Value = 28.02000045776367
$ws.Range("H37").Value = 112123800
$ws.Range("I37").Value = "TSEM"

$ws.Range("D38").Value = 33.65000152587891
$ws.Range("E38").Value = 32.86999893188477
$ws.Range("F38").Value = 34.59000015258789
$ws.Range("G38").Value = 30.73999977111816
$ws.Range("H38").Value = 112123800
$ws.Range("I38").Value = "TSEM"

$ws.Range("D39").Value = 39.31000137329102
$ws.Range("E39").Value = 40.79999923706055
$ws.Range("F39").Value = 44.18000030517578
$ws.Range("G39").Value = 37.88999938964844
$ws.Range("H39").Value = 112123800
$ws.Range("I39").Value = "TSEM"

$ws.Range("D40").Value = 44.09000015258789
$ws.Range("E40").Value = 41.95999908447266
$ws.Range("F40").Value = 46.54000091552734
$ws.Range("G40").Value = 41.54000091552734
$ws.Range("H40").Value = 112123800
$ws.Range("I40").Value = "TSEM"

$ws.Range("D41").Value = 52
$ws.Range("E41").Value = 48.88999938964844
$ws.Range("F41").Value = 55.31000137329102
$ws.Range("G41").Value = 44.52000045776367
$ws.Range("H41").Value = 112123800
$ws.Range("I41").Value = "TSEM"

$ws.Range("D42").Value = 35.54999923706055
$ws.Range("E42").Value = 35.77999877929688
$ws.Range("F42").Value = 37.77999877929688
$ws.Range("G42").Value = 28.63999938964844
$ws.Range("H42").Value = 112123800
$ws.Range("I42").Value = "TSEM"

$ws.Range("D43").Value = 43.36000061035156
$ws.Range("E43").Value = 45.75
$ws.Range("F43").Value = 50.93000030517578
$ws.Range("G43").Value = 42.08000183105469
$ws.Range("H43").Value = 112123800
$ws.Range("I43").Value = "TSEM"
